$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new text value. D/E columns are plain-text price/
# percentage strings in the source sheet (t="inlineStr"), so force the
# range to Text format first -- otherwise Excel's COM Value setter will
# auto-coerce numeric-looking strings (e.g. "1.00", "579.03") into real
# numbers and mangle/round them ("1", "579.0299999999...").
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.945.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.456.85'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.63%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.03'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.93'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.25%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.456.82'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.88'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.331'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.901.99'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.28'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.56%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.779.27'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.453.38'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.69'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.03'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.09%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.31'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -7.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.75'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -7.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.85'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -10.11%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.582.27'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0894'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -7.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '502.74'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -5.29%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -5.46%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.54%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '158.99'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -9.17%  '
$ws.Range("B38").Value = 'WhiteBITCoin'
$ws.Range("C38").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.56'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.42'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.33'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.24%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.325'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -6.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.75'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '38.80'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.30'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -7.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.06'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.72%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.15%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0732'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.03%  '
